$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.934.32'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.553.86'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").Value = '''206.97'
$ws.Range("D5").Style = $ws.Range("A1").Style
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '''21.76'
$ws.Range("D8").Style = $ws.Range("A1").Style
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '''0.0586'
$ws.Range("D10").Style = $ws.Range("A1").Style
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").Value = '''0.0859'
$ws.Range("D11").Style = $ws.Range("A1").Style
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").Value = '1.775.01'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").Value = '1.557.41'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '''0.517'
$ws.Range("D15").Style = $ws.Range("A1").Style
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").Value = '26.937.66'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '''61.69'
$ws.Range("D17").Style = $ws.Range("A1").Style
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '''217.09'
$ws.Range("D18").Style = $ws.Range("A1").Style
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("D19").Value = '0.0₃0688'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '''7.23'
$ws.Range("D20").Style = $ws.Range("A1").Style
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("E22").Value = '  +1.57%  '
$ws.Range("D23").Value = '''9.19'
$ws.Range("D23").Style = $ws.Range("A1").Style
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").Value = '''153.82'
$ws.Range("D25").Style = $ws.Range("A1").Style
$ws.Range("E25").Value = '  +1.64%  '
$ws.Range("D26").Value = '''6.57'
$ws.Range("D26").Style = $ws.Range("A1").Style
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '''14.87'
$ws.Range("D27").Style = $ws.Range("A1").Style
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").Value = '''0.0468'
$ws.Range("D30").Style = $ws.Range("A1").Style
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").Value = '1.423.67'
$ws.Range("E33").Value = '  +4.59%  '
$ws.Range("D34").Value = '''3.02'
$ws.Range("D34").Style = $ws.Range("A1").Style
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("E35").Value = '  +3.83%  '
$ws.Range("D36").Value = '''0.961'
$ws.Range("D36").Style = $ws.Range("A1").Style
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").Value = '''0.520'
$ws.Range("D39").Style = $ws.Range("A1").Style
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  +1.31%  '
$ws.Range("D42").Value = '''5.65'
$ws.Range("D42").Style = $ws.Range("A1").Style
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").Value = '''0.987'
$ws.Range("D43").Style = $ws.Range("A1").Style
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("E44").Value = '  +3.74%  '
$ws.Range("D45").Value = '''63.80'
$ws.Range("D45").Style = $ws.Range("A1").Style
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = '1.689.71'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '''0.0523'
$ws.Range("D49").Style = $ws.Range("A1").Style
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("D50").Value = '0.0₆0101'
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("E51").Value = '  +1.58%  '
